$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (Taxa) for rows 2-4: 0 -> 2
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 2

# Update column C (Data de Salvamento) for rows 2-15: seconds 25 -> 35
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = "2025-04-04 13:26:35"
}
